$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 445
$ws.Range("I2").Value = 445
$ws.Range("K2").Value = 445
$ws.Range("M2").Value = -332
$ws.Range("H10").Value = 5250
$ws.Range("J10").Value = 5250
$ws.Range("L10").Value = 5250
$ws.Range("N10").Value = -5836
$ws.Range("H11").Value = 109.333336
$ws.Range("I11").Value = 109.333336
$ws.Range("K11").Value = 109.333336
$ws.Range("M11").Value = 30.666664
$ws.Range("H13").Value = 2500
$ws.Range("J13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("N13").Value = -2838
$ws.Range("H15").Value = 203.9375
$ws.Range("I15").Value = 203.9375
$ws.Range("K15").Value = 611.8125
$ws.Range("M15").Value = -442.8125
$ws.Range("H29").Value = 2400
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 3500
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 10500
$ws.Range("M29").Value = -319
$ws.Range("N29").Value = -11062
$ws.Range("H74").Value = 3749.75
$ws.Range("I74").Value = 3499.5
$ws.Range("K74").Value = 3499.5
$ws.Range("M74").Value = -2563.5
$ws.Range("H77").Value = 3749.75
$ws.Range("I77").Value = 3499.5
$ws.Range("K77").Value = 17497.5
$ws.Range("M77").Value = -12817.5
$ws.Range("H80").Value = 809
$ws.Range("I80").Value = 583.1667
$ws.Range("J80").Value = 1080
$ws.Range("K80").Value = 1749.5001
$ws.Range("L80").Value = 3240
$ws.Range("M80").Value = -751.5001
$ws.Range("N80").Value = -5236
$ws.Range("H83").Value = 809
$ws.Range("I83").Value = 583.1667
$ws.Range("J83").Value = 1080
$ws.Range("K83").Value = 5248.5003
$ws.Range("L83").Value = 9720
$ws.Range("M83").Value = -256.5002999999997
$ws.Range("N83").Value = -19704
$ws.Range("H100").Value = 1533.8889
$ws.Range("I100").Value = 1724.875
$ws.Range("J100").Value = 6
$ws.Range("K100").Value = 1724.875
$ws.Range("L100").Value = 6
$ws.Range("M100").Value = -1183.875
$ws.Range("N100").Value = -1088
$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988
$ws.Range("H107").Value = 1598.9
$ws.Range("I107").Value = 936.25
$ws.Range("K107").Value = 936.25
$ws.Range("M107").Value = 983.75
$ws.Range("H127").Value = 738.5
$ws.Range("I127").Value = 670.1429000000001
$ws.Range("K127").Value = 2010.4287
$ws.Range("M127").Value = 2949.5713
$ws.Range("H129").Value = 1183
$ws.Range("I129").Value = 1183
$ws.Range("K129").Value = 3549
$ws.Range("M129").Value = 1451
$ws.Range("H132").Value = 4375.615
$ws.Range("I132").Value = 4530.091
$ws.Range("K132").Value = 13590.273
$ws.Range("M132").Value = -11060.273
$ws.Range("H135").Value = 1376.9412
$ws.Range("I135").Value = 882.0833
$ws.Range("K135").Value = 7938.7497
$ws.Range("M135").Value = -5403.7497
$ws.Range("H137").Value = 3287
$ws.Range("I137").Value = 2259.9
$ws.Range("J137").Value = 3891.1765
$ws.Range("K137").Value = 6779.700000000001
$ws.Range("L137").Value = 11673.5295
$ws.Range("M137").Value = -4229.700000000001
$ws.Range("N137").Value = -16773.5295
$ws.Range("H138").Value = 4599.263
$ws.Range("J138").Value = 4969.091
$ws.Range("L138").Value = 14907.273
$ws.Range("N138").Value = -25187.273

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6267.706
$ws.Range("I32").Value = 4437.6665
$ws.Range("J32").Value = 19993
$ws.Range("K32").Value = 4437.6665
$ws.Range("L32").Value = 19993
$ws.Range("M32").Value = -4150.6665
$ws.Range("N32").Value = -20567
$ws.Range("H97").Value = 942.5
$ws.Range("I97").Value = 573.63635
$ws.Range("K97").Value = 573.63635
$ws.Range("M97").Value = -77.63634999999999

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1206.8572
$ws.Range("I64").Value = 1616
$ws.Range("J64").Value = 900
$ws.Range("K64").Value = 1616
$ws.Range("L64").Value = 900
$ws.Range("M64").Value = -1391
$ws.Range("N64").Value = -1350
$ws.Range("H67").Value = 1206.8572
$ws.Range("I67").Value = 1616
$ws.Range("J67").Value = 900
$ws.Range("K67").Value = 1616
$ws.Range("L67").Value = 900
$ws.Range("M67").Value = -836
$ws.Range("N67").Value = -2460
$ws.Range("H80").Value = 873.8
$ws.Range("I80").Value = 788
$ws.Range("K80").Value = 788
$ws.Range("M80").Value = 210
$ws.Range("H83").Value = 873.8
$ws.Range("I83").Value = 788
$ws.Range("K83").Value = 3940
$ws.Range("M83").Value = 1052

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H62").Value = 2965
$ws.Range("I62").Value = 2950
$ws.Range("K62").Value = 2950
$ws.Range("M62").Value = -2326
$ws.Range("H65").Value = 2965
$ws.Range("I65").Value = 2950
$ws.Range("K65").Value = 14750
$ws.Range("M65").Value = -11630
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H134").Value = 1463
$ws.Range("I134").Value = 1463
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4389
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -1854

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 16999.5
$ws.Range("I120").Value = 16999.5
$ws.Range("K120").Value = 50998.5
$ws.Range("M120").Value = -46160.5
$ws.Range("H129").Value = 1038
$ws.Range("I129").Value = 1307.5
$ws.Range("J129").Value = 499
$ws.Range("K129").Value = 3922.5
$ws.Range("L129").Value = 1497
$ws.Range("M129").Value = 1077.5
$ws.Range("N129").Value = -11497
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").ClearContents()
$ws.Range("N132").Value = 0

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5212.143
$ws.Range("I102").Value = 5228.4614
$ws.Range("K102").Value = 5228.4614
$ws.Range("M102").Value = -3606.4614
$ws.Range("H126").Value = 4765.6665
$ws.Range("I126").Value = 4979.5
$ws.Range("J126").Value = 4338
$ws.Range("K126").Value = 14938.5
$ws.Range("L126").Value = 13014
$ws.Range("M126").Value = -12468.5
$ws.Range("N126").Value = -17954

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5259
$ws.Range("I7").Value = 5259
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 5259
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -5147
$ws.Range("H16").Value = 428.77777
$ws.Range("I16").Value = 448.625
$ws.Range("J16").Value = 270
$ws.Range("K16").Value = 448.625
$ws.Range("L16").Value = 270
$ws.Range("M16").Value = -278.625
$ws.Range("N16").Value = -610
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H55").Value = 3334
$ws.Range("I55").Value = 3810
$ws.Range("J55").Value = 2
$ws.Range("K55").Value = 3810
$ws.Range("L55").Value = 2
$ws.Range("M55").Value = -3637
$ws.Range("N55").Value = -348
$ws.Range("H68").Value = 3012.125
$ws.Range("I68").Value = 2942.4285
$ws.Range("J68").Value = 3500
$ws.Range("K68").Value = 2942.4285
$ws.Range("L68").Value = 3500
$ws.Range("M68").Value = -2193.4285
$ws.Range("N68").Value = -4998
$ws.Range("H71").Value = 3012.125
$ws.Range("I71").Value = 2942.4285
$ws.Range("J71").Value = 3500
$ws.Range("K71").Value = 14712.1425
$ws.Range("L71").Value = 17500
$ws.Range("M71").Value = -10968.1425
$ws.Range("N71").Value = -24988
$ws.Range("H82").Value = 1964.7
$ws.Range("J82").Value = 2500
$ws.Range("L82").Value = 2500
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 1964.7
$ws.Range("J85").Value = 2500
$ws.Range("L85").Value = 2500
$ws.Range("N85").Value = -4996
$ws.Range("H100").Value = 4200
$ws.Range("I100").Value = 3875
$ws.Range("K100").Value = 3875
$ws.Range("M100").Value = -3334
$ws.Range("H126").Value = 5259
$ws.Range("I126").Value = 5259
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 15777
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13307
$ws.Range("H132").Value = 36666.5
$ws.Range("I132").Value = 33999.8
$ws.Range("J132").Value = 50000
$ws.Range("K132").Value = 101999.4
$ws.Range("L132").Value = 150000
$ws.Range("M132").Value = -99469.40000000001
$ws.Range("N132").Value = -155060

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3500.6667
$ws.Range("I62").Value = 3251
$ws.Range("K62").Value = 3251
$ws.Range("M62").Value = -2627
$ws.Range("H65").Value = 3500.6667
$ws.Range("I65").Value = 3251
$ws.Range("K65").Value = 16255
$ws.Range("M65").Value = -13135
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H136").Value = 1815.6666
$ws.Range("I136").Value = 1815.6666
$ws.Range("K136").Value = 5446.9998
$ws.Range("M136").Value = -2896.9998
